$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 31.14997866666667
$ws.Range("H2").Value = 93.44993600000001
$ws.Range("I2").Value = 0.4621739036316256
$ws.Range("J2").Value = 0.4621739036316256
$ws.Range("M2").Value = 0.7004666666666667
$ws.Range("N2").Value = 2.1014
$ws.Range("O2").Value = 0.08010868614664106
$ws.Range("P2").Value = 0.08010868614664106
$ws.Range("Q2").Value = 21.81952172337778
$ws.Range("R2").Value = 196.3756955104
$ws.Range("S2").Value = 0.03702414419119383
$ws.Range("T2").Value = 0.03702414419119383
$ws.Range("G3").Value = 31.14997866666667
$ws.Range("H3").Value = 93.44993600000001
$ws.Range("I3").Value = 0.4621739036316256
$ws.Range("J3").Value = 0.4621739036316256
$ws.Range("M3").Value = 3.474244666666667
$ws.Range("N3").Value = 10.422734
$ws.Range("O3").Value = 0.3973310777557461
$ws.Range("P3").Value = 0.3973310777557461
$ws.Range("Q3").Value = 108.2226472494471
$ws.Range("R3").Value = 974.0038252450241
$ws.Range("S3").Value = 0.1836360552405341
$ws.Range("T3").Value = 0.1836360552405341
$ws.Range("G4").Value = 31.14997866666667
$ws.Range("H4").Value = 93.44993600000001
$ws.Range("I4").Value = 0.4621739036316256
$ws.Range("J4").Value = 0.4621739036316256
$ws.Range("M4").Value = 2.336226333333333
$ws.Range("N4").Value = 7.008679
$ws.Range("O4").Value = 0.2671819103043467
$ws.Range("P4").Value = 0.2671819103043467
$ws.Range("Q4").Value = 72.77340044383823
$ws.Range("R4").Value = 654.960603994544
$ws.Range("S4").Value = 0.1234845064651148
$ws.Range("T4").Value = 0.1234845064651148
$ws.Range("G5").Value = 31.14997866666667
$ws.Range("H5").Value = 93.44993600000001
$ws.Range("I5").Value = 0.4621739036316256
$ws.Range("J5").Value = 0.4621739036316256
$ws.Range("M5").Value = 2.233016333333333
$ws.Range("N5").Value = 6.699049
$ws.Range("O5").Value = 0.2553783257932662
$ws.Range("P5").Value = 0.2553783257932662
$ws.Range("Q5").Value = 69.55841114565156
$ws.Range("R5").Value = 626.0257003108641
$ws.Range("S5").Value = 0.1180291977347829
$ws.Range("T5").Value = 0.1180291977347829
$ws.Range("G6").Value = 18.94069966666667
$ws.Range("H6").Value = 56.822099
$ws.Range("I6").Value = 0.2810241764892454
$ws.Range("J6").Value = 0.2810241764892454
$ws.Range("M6").Value = 0.7004666666666667
$ws.Range("N6").Value = 2.1014
$ws.Range("O6").Value = 0.08010868614664106
$ws.Range("P6").Value = 0.08010868614664106
$ws.Range("Q6").Value = 13.26732875984444
$ws.Range("R6").Value = 119.4059588386
$ws.Range("S6").Value = 0.02251247755399522
$ws.Range("T6").Value = 0.02251247755399522
$ws.Range("G7").Value = 18.94069966666667
$ws.Range("H7").Value = 56.822099
$ws.Range("I7").Value = 0.2810241764892454
$ws.Range("J7").Value = 0.2810241764892454
$ws.Range("M7").Value = 3.474244666666667
$ws.Range("N7").Value = 10.422734
$ws.Range("O7").Value = 0.3973310777557461
$ws.Range("P7").Value = 0.3973310777557461
$ws.Range("Q7").Value = 65.80462479985178
$ws.Range("R7").Value = 592.241623198666
$ws.Range("S7").Value = 0.1116596389198929
$ws.Range("T7").Value = 0.1116596389198929
$ws.Range("G8").Value = 18.94069966666667
$ws.Range("H8").Value = 56.822099
$ws.Range("I8").Value = 0.2810241764892454
$ws.Range("J8").Value = 0.2810241764892454
$ws.Range("M8").Value = 2.336226333333333
$ws.Range("N8").Value = 7.008679
$ws.Range("O8").Value = 0.2671819103043467
$ws.Range("P8").Value = 0.2671819103043467
$ws.Range("Q8").Value = 44.24976133302456
$ws.Range("R8").Value = 398.247851997221
$ws.Range("S8").Value = 0.07508457631610245
$ws.Range("T8").Value = 0.07508457631610245
$ws.Range("G9").Value = 18.94069966666667
$ws.Range("H9").Value = 56.822099
$ws.Range("I9").Value = 0.2810241764892454
$ws.Range("J9").Value = 0.2810241764892454
$ws.Range("M9").Value = 2.233016333333333
$ws.Range("N9").Value = 6.699049
$ws.Range("O9").Value = 0.2553783257932662
$ws.Range("P9").Value = 0.2553783257932662
$ws.Range("Q9").Value = 42.29489172042789
$ws.Range("R9").Value = 380.6540254838511
$ws.Range("S9").Value = 0.07176748369925486
$ws.Range("T9").Value = 0.07176748369925486
$ws.Range("G10").Value = 14.86848
$ws.Range("H10").Value = 44.60544
$ws.Range("I10").Value = 0.2206044349565553
$ws.Range("J10").Value = 0.2206044349565553
$ws.Range("M10").Value = 0.7004666666666667
$ws.Range("N10").Value = 2.1014
$ws.Range("O10").Value = 0.08010868614664106
$ws.Range("P10").Value = 0.08010868614664106
$ws.Range("Q10").Value = 10.414874624
$ws.Range("R10").Value = 93.733871616
$ws.Range("S10").Value = 0.01767233144249178
$ws.Range("T10").Value = 0.01767233144249178
$ws.Range("G11").Value = 14.86848
$ws.Range("H11").Value = 44.60544
$ws.Range("I11").Value = 0.2206044349565553
$ws.Range("J11").Value = 0.2206044349565553
$ws.Range("M11").Value = 3.474244666666667
$ws.Range("N11").Value = 10.422734
$ws.Range("O11").Value = 0.3973310777557461
$ws.Range("P11").Value = 0.3973310777557461
$ws.Range("Q11").Value = 51.65673734144
$ws.Range("R11").Value = 464.91063607296
$ws.Range("S11").Value = 0.0876529978989855
$ws.Range("T11").Value = 0.0876529978989855
$ws.Range("G12").Value = 14.86848
$ws.Range("H12").Value = 44.60544
$ws.Range("I12").Value = 0.2206044349565553
$ws.Range("J12").Value = 0.2206044349565553
$ws.Range("M12").Value = 2.336226333333333
$ws.Range("N12").Value = 7.008679
$ws.Range("O12").Value = 0.2671819103043467
$ws.Range("P12").Value = 0.2671819103043467
$ws.Range("Q12").Value = 34.73613451264
$ws.Range("R12").Value = 312.62521061376
$ws.Range("S12").Value = 0.05894151435330344
$ws.Range("T12").Value = 0.05894151435330344
$ws.Range("G13").Value = 14.86848
$ws.Range("H13").Value = 44.60544
$ws.Range("I13").Value = 0.2206044349565553
$ws.Range("J13").Value = 0.2206044349565553
$ws.Range("M13").Value = 2.233016333333333
$ws.Range("N13").Value = 6.699049
$ws.Range("O13").Value = 0.2553783257932662
$ws.Range("P13").Value = 0.2553783257932662
$ws.Range("Q13").Value = 33.20155869184
$ws.Range("R13").Value = 298.81402822656
$ws.Range("S13").Value = 0.05633759126177459
$ws.Range("T13").Value = 0.05633759126177459
$ws.Range("G14").Value = 2.439668
$ws.Range("H14").Value = 7.319004
$ws.Range("I14").Value = 0.03619748492257375
$ws.Range("J14").Value = 0.03619748492257375
$ws.Range("M14").Value = 0.7004666666666667
$ws.Range("N14").Value = 2.1014
$ws.Range("O14").Value = 0.08010868614664106
$ws.Range("P14").Value = 0.08010868614664106
$ws.Range("Q14").Value = 1.708906111733333
$ws.Range("R14").Value = 15.3801550056
$ws.Range("S14").Value = 0.002899732958960232
$ws.Range("T14").Value = 0.002899732958960232
$ws.Range("G15").Value = 2.439668
$ws.Range("H15").Value = 7.319004
$ws.Range("I15").Value = 0.03619748492257375
$ws.Range("J15").Value = 0.03619748492257375
$ws.Range("M15").Value = 3.474244666666667
$ws.Range("N15").Value = 10.422734
$ws.Range("O15").Value = 0.3973310777557461
$ws.Range("P15").Value = 0.3973310777557461
$ws.Range("Q15").Value = 8.476003537437332
$ws.Range("R15").Value = 76.284031836936
$ws.Range("S15").Value = 0.0143823856963336
$ws.Range("T15").Value = 0.0143823856963336
$ws.Range("G16").Value = 2.439668
$ws.Range("H16").Value = 7.319004
$ws.Range("I16").Value = 0.03619748492257375
$ws.Range("J16").Value = 0.03619748492257375
$ws.Range("M16").Value = 2.336226333333333
$ws.Range("N16").Value = 7.008679
$ws.Range("O16").Value = 0.2671819103043467
$ws.Range("P16").Value = 0.2671819103043467
$ws.Range("Q16").Value = 5.699616626190666
$ws.Range("R16").Value = 51.296549635716
$ws.Range("S16").Value = 0.009671313169826041
$ws.Range("T16").Value = 0.009671313169826041
$ws.Range("G17").Value = 2.439668
$ws.Range("H17").Value = 7.319004
$ws.Range("I17").Value = 0.03619748492257375
$ws.Range("J17").Value = 0.03619748492257375
$ws.Range("M17").Value = 2.233016333333333
$ws.Range("N17").Value = 6.699049
$ws.Range("O17").Value = 0.2553783257932662
$ws.Range("P17").Value = 0.2553783257932662
$ws.Range("Q17").Value = 5.447818491910667
$ws.Range("R17").Value = 49.030366427196
$ws.Range("S17").Value = 0.009244053097453882
$ws.Range("T17").Value = 0.009244053097453882
